$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name labels (column A) whose rank/position changed ---
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 24 de Marzo de 2020 a las 18:46'
$ws.Cells.Item(18, 1).Value = 'Canada'
$ws.Cells.Item(19, 1).Value = 'Portugal'
$ws.Cells.Item(20, 1).Value = 'Suecia'
$ws.Cells.Item(69, 1).Value = 'Lituania'
$ws.Cells.Item(70, 1).Value = 'Eslovaquia'
$ws.Cells.Item(75, 1).Value = 'Marruecos'
$ws.Cells.Item(76, 1).Value = 'Principado de Andorra'
$ws.Cells.Item(77, 1).Value = 'Bosnia y Herzegovina'
$ws.Cells.Item(78, 1).Value = 'Uruguay'
$ws.Cells.Item(79, 1).Value = 'Costa Rica'
$ws.Cells.Item(80, 1).Value = 'Nueva Zelanda'
$ws.Cells.Item(81, 1).Value = 'Republica de Macedonia'
$ws.Cells.Item(113, 1).Value = 'Montenegro'
$ws.Cells.Item(114, 1).Value = 'Consejo Danes para los Refugiados'
$ws.Cells.Item(115, 1).Value = 'Kirguistan'
$ws.Cells.Item(116, 1).Value = 'Mauricio'
$ws.Cells.Item(117, 1).Value = 'Nigeria'
$ws.Cells.Item(118, 1).Value = 'Puerto Rico'
$ws.Cells.Item(119, 1).Value = 'Banglades'
$ws.Cells.Item(120, 1).Value = 'Ruanda'
$ws.Cells.Item(121, 1).Value = 'Mayotte'
$ws.Cells.Item(122, 1).Value = 'Guam'
$ws.Cells.Item(123, 1).Value = 'Honduras'
$ws.Cells.Item(125, 1).Value = 'Macao'
$ws.Cells.Item(126, 1).Value = 'Paraguay'
$ws.Cells.Item(129, 1).Value = 'Polinesia Francesa'
$ws.Cells.Item(130, 1).Value = 'Isla de Man'
$ws.Cells.Item(137, 1).Value = 'Islas Virgenes de los Estados Unidos'
$ws.Cells.Item(138, 1).Value = 'Madagascar'
$ws.Cells.Item(139, 1).Value = 'Barbados'
$ws.Cells.Item(142, 1).Value = 'Tanzania'
$ws.Cells.Item(143, 1).Value = 'Etiopia'
$ws.Cells.Item(147, 1).Value = 'Uganda'
$ws.Cells.Item(148, 1).Value = 'Guinea Ecuatorial'
$ws.Cells.Item(151, 1).Value = 'Benin'
$ws.Cells.Item(152, 1).Value = 'Haiti'
$ws.Cells.Item(153, 1).Value = 'Bermudas'
$ws.Cells.Item(154, 1).Value = 'Surinam'
$ws.Cells.Item(155, 1).Value = 'Gabon'
$ws.Cells.Item(156, 1).Value = 'Namibia'
$ws.Cells.Item(161, 1).Value = 'Suazilandia'
$ws.Cells.Item(162, 1).Value = 'Bahamas'

# --- Update statistic values (columns B-H) ---
$ws.Cells.Item(8, 2).Value = 32781
$ws.Cells.Item(8, 3).Value = 3725
$ws.Cells.Item(8, 4).Value = 3133
$ws.Cells.Item(8, 5).Value = 29492
$ws.Cells.Item(8, 7).Value = 33
$ws.Cells.Item(8, 8).Value = 156
$ws.Cells.Item(18, 2).Value = 2583
$ws.Cells.Item(18, 3).Value = 492
$ws.Cells.Item(18, 4).Value = 112
$ws.Cells.Item(18, 5).Value = 2447
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 24
$ws.Cells.Item(19, 2).Value = 2362
$ws.Cells.Item(19, 3).Value = 302
$ws.Cells.Item(19, 4).Value = 22
$ws.Cells.Item(19, 5).Value = 2307
$ws.Cells.Item(19, 6).Value = 48
$ws.Cells.Item(19, 7).Value = 10
$ws.Cells.Item(19, 8).Value = 33
$ws.Cells.Item(20, 2).Value = 2286
$ws.Cells.Item(20, 3).Value = 240
$ws.Cells.Item(20, 4).Value = 16
$ws.Cells.Item(20, 5).Value = 2234
$ws.Cells.Item(20, 6).Value = 136
$ws.Cells.Item(20, 7).Value = 9
$ws.Cells.Item(20, 8).Value = 36
$ws.Cells.Item(23, 5).Value = 1604
$ws.Cells.Item(23, 7).Value = 2
$ws.Cells.Item(23, 8).Value = 3
$ws.Cells.Item(51, 4).Value = 177
$ws.Cells.Item(51, 5).Value = 210
$ws.Cells.Item(64, 4).Value = 24
$ws.Cells.Item(64, 5).Value = 221
$ws.Cells.Item(69, 2).Value = 209
$ws.Cells.Item(69, 3).Value = 30
$ws.Cells.Item(69, 4).Value = 1
$ws.Cells.Item(69, 5).Value = 207
$ws.Cells.Item(69, 6).Value = 1
$ws.Cells.Item(69, 8).Value = 1
$ws.Cells.Item(70, 2).Value = 204
$ws.Cells.Item(70, 3).Value = 18
$ws.Cells.Item(70, 4).Value = 7
$ws.Cells.Item(70, 5).Value = 197
$ws.Cells.Item(70, 6).Value = 2
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(75, 2).Value = 170
$ws.Cells.Item(75, 3).Value = 27
$ws.Cells.Item(75, 4).Value = 6
$ws.Cells.Item(75, 5).Value = 159
$ws.Cells.Item(75, 6).Value = 1
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = 5
$ws.Cells.Item(76, 2).Value = 164
$ws.Cells.Item(76, 3).Value = 31
$ws.Cells.Item(76, 4).Value = 1
$ws.Cells.Item(76, 6).Value = 7
$ws.Cells.Item(76, 8).Value = 1
$ws.Cells.Item(77, 2).Value = 164
$ws.Cells.Item(77, 3).Value = 28
$ws.Cells.Item(77, 5).Value = 160
$ws.Cells.Item(77, 6).Value = 1
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(78, 2).Value = 162
$ws.Cells.Item(78, 3).Value = 0
$ws.Cells.Item(78, 4).Value = 0
$ws.Cells.Item(78, 5).Value = 162
$ws.Cells.Item(78, 6).Value = 3
$ws.Cells.Item(79, 2).Value = 158
$ws.Cells.Item(79, 3).Value = 0
$ws.Cells.Item(79, 5).Value = 154
$ws.Cells.Item(79, 6).Value = 2
$ws.Cells.Item(79, 7).Value = 0
$ws.Cells.Item(80, 2).Value = 155
$ws.Cells.Item(80, 3).Value = 53
$ws.Cells.Item(80, 4).Value = 12
$ws.Cells.Item(80, 5).Value = 143
$ws.Cells.Item(80, 6).Value = 0
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(81, 2).Value = 148
$ws.Cells.Item(81, 3).Value = 12
$ws.Cells.Item(81, 4).Value = 1
$ws.Cells.Item(81, 5).Value = 145
$ws.Cells.Item(81, 8).Value = 2
$ws.Cells.Item(113, 2).Value = 47
$ws.Cells.Item(113, 3).Value = 20
$ws.Cells.Item(113, 5).Value = 46
$ws.Cells.Item(113, 8).Value = 1
$ws.Cells.Item(114, 2).Value = 45
$ws.Cells.Item(114, 3).Value = 9
$ws.Cells.Item(114, 5).Value = 43
$ws.Cells.Item(114, 8).Value = 2
$ws.Cells.Item(115, 3).Value = 26
$ws.Cells.Item(115, 5).Value = 42
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 8).Value = 0
$ws.Cells.Item(116, 3).Value = 6
$ws.Cells.Item(116, 4).Value = 0
$ws.Cells.Item(116, 5).Value = 40
$ws.Cells.Item(116, 6).Value = 1
$ws.Cells.Item(116, 8).Value = 2
$ws.Cells.Item(117, 2).Value = 42
$ws.Cells.Item(117, 3).Value = 2
$ws.Cells.Item(117, 4).Value = 2
$ws.Cells.Item(117, 5).Value = 39
$ws.Cells.Item(117, 8).Value = 1
$ws.Cells.Item(118, 3).Value = 8
$ws.Cells.Item(118, 4).Value = 1
$ws.Cells.Item(118, 5).Value = 36
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 2
$ws.Cells.Item(119, 2).Value = 39
$ws.Cells.Item(119, 3).Value = 6
$ws.Cells.Item(119, 4).Value = 5
$ws.Cells.Item(119, 5).Value = 30
$ws.Cells.Item(119, 7).Value = 1
$ws.Cells.Item(119, 8).Value = 4
$ws.Cells.Item(120, 3).Value = 0
$ws.Cells.Item(121, 2).Value = 36
$ws.Cells.Item(121, 3).Value = 12
$ws.Cells.Item(121, 5).Value = 36
$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(122, 2).Value = 32
$ws.Cells.Item(122, 3).Value = 3
$ws.Cells.Item(122, 5).Value = 31
$ws.Cells.Item(122, 8).Value = 1
$ws.Cells.Item(123, 2).Value = 30
$ws.Cells.Item(123, 3).Value = 0
$ws.Cells.Item(123, 5).Value = 30
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(125, 2).Value = 28
$ws.Cells.Item(125, 3).Value = 3
$ws.Cells.Item(125, 4).Value = 10
$ws.Cells.Item(125, 5).Value = 18
$ws.Cells.Item(125, 6).Value = 0
$ws.Cells.Item(125, 7).Value = 0
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(126, 2).Value = 27
$ws.Cells.Item(126, 3).Value = 5
$ws.Cells.Item(126, 4).Value = 0
$ws.Cells.Item(126, 5).Value = 25
$ws.Cells.Item(126, 6).Value = 1
$ws.Cells.Item(126, 7).Value = 1
$ws.Cells.Item(126, 8).Value = 2
$ws.Cells.Item(129, 3).Value = 5
$ws.Cells.Item(130, 3).Value = 10
$ws.Cells.Item(137, 3).Value = 0
$ws.Cells.Item(138, 3).Value = 5
$ws.Cells.Item(142, 3).Value = 0
$ws.Cells.Item(143, 3).Value = 1
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(154, 3).Value = 1
$ws.Cells.Item(155, 5).Value = 5
$ws.Cells.Item(155, 8).Value = 1
$ws.Cells.Item(156, 3).Value = 2
$ws.Cells.Item(156, 4).Value = 2
$ws.Cells.Item(156, 5).Value = 4
$ws.Cells.Item(156, 8).Value = 0
